$wb = $excel.ActiveWorkbook

# Add a "Comments" header column (E1) to the history sheets, selecting E1
# on each (matching the user clicking/typing into the new header cell).
$sheetNames = @("Withdraw History", "Deposit History", "Transfer History")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Activate()
    $ws.Range("E1").Value = "Comments"
    $ws.Range("E1").Select() | Out-Null
}

# The active sheet (Absolute History) ends up with the selection one cell
# below the newly entered header, matching interactive entry (Enter moves
# down after typing).
$activeWs = $wb.Worksheets.Item("Absolute History")
$activeWs.Activate()
$activeWs.Range("E1").Value = "Comments"
$activeWs.Range("E2").Select() | Out-Null
